$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new year column S (2022), copying the header/data formatting from column R
$ws.Range("R4").Copy()
$ws.Range("S4").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("R5").Copy()
$ws.Range("S5").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Set the new header value (2022) and new data value (3.4)
$ws.Range("S4").Value = 2022
$ws.Range("S5").Value = 3.4

# Refresh the last three years of data with the updated figures
$ws.Range("P5").Value = 4.4000000000000004
$ws.Range("Q5").Value = 2.9
$ws.Range("R5").Value = 3.2

# Update the active selection to match the new extent
$ws.Range("T4").Select()
